$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new weekly data row at row 45 -- this shifts the existing
# rows 45..109 down to 46..110 (preserving all of their data/formatting)
# and leaves a blank row 45 (inheriting the date number-format from the
# row above) ready to be populated with the new "Orégano" quote.
$ws.Rows.Item(45).Insert()

$ws.Range("A45").Value = 6
$ws.Range("B45").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C45").Value = "Metropolitana"
$ws.Range("D45").Value = 44467
$ws.Range("E45").Value = 13
$ws.Range("F45").Value = 100112029
$ws.Range("G45").Value = "Orégano"
$ws.Range("H45").Value = "Sin especificar"
$ws.Range("I45").Value = "Primera"
$ws.Range("J45").Value = 33
$ws.Range("K45").Value = 9000
$ws.Range("L45").Value = 10000
$ws.Range("M45").Value = 9515
$ws.Range("N45").Value = "$/docena de atados"
$ws.Range("O45").Value = "Región Metropolitana"
$ws.Range("P45").Value = 3172
$ws.Range("Q45").Value = 3
$ws.Range("R45").Value = "Hortaliza"
